$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the two missing point values (rows 6 and 7 of column B)
$ws.Range("B6").Value = 17
$ws.Range("B7").Value = 20

# Row 11 keeps the old (pre-resize) row height explicitly
$ws.Rows.Item(11).RowHeight = 15

# Move the active selection to E5 (as recorded in the saved view state)
$ws.Range("E5").Select() | Out-Null
